$d = $word.ActiveDocument

# --- Edit 1: remove gramStart/gramEnd proofErr markers around
# "If needed we will also work in weekends." ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*If needed we will also work in weekends.*") {
        $fixedXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="001879D9" w:rsidRPr="00DB6FDD" w:rsidRDefault="00DB6FDD"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00DB6FDD"><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>If needed we will also work in weekends.</w:t></w:r></w:p>'
        $p.Range.InsertXML($fixedXml)
        break
    }
}

# --- Edit 2: append the new "Log 2013-04-25" section at the end of the
# document body (before the final section break). ---
$insertXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Log 2013-04-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>25</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Driver: Simon</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Navigator: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Thelle</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">What did we </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>do:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Placed objects in the world</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Added exit descriptions to location</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Improved capturing user input</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Plan for next time:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Make it possible for the player to pick up and use items</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>JUnit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> test of added features</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertXML($insertXml)

Write-Output "edit complete"
